$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")

# Mark the empty make/id cells on row 2 with a placeholder asterisk
$ws.Range("A2").Value = "*"
$ws.Range("B2").Value = "*"

# Move the active selection from C3 to B3
$ws.Range("B3").Select()

